$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# TC02 query cell: A2 now carries the Neo4j MATCH query that drives this test case
$query = 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.gender IN [''FEMALE''] RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(t.clinical_trial_designation ,'''')as `Trial Code` , coalesce(a.arm_id,'''') As `Arm` , coalesce(a.arm_drug,'''') As `Arm Treatment` , coalesce(c.disease,'''') As Diagnosis , coalesce(c.gender,'''') As Gender , coalesce(c.race,'''') As Race , coalesce(c.ethnicity,'''') As Ethnicity'

$ws.Range("A2").Value = $query
$ws.Range("A2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 87

# Selection moves onto the newly-populated query column
$ws.Range("A2:A5").Select()
